$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '41.525.79'
$ws.Range('E2').Value = '  -1.07%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.163.79'
$ws.Range('E3').Value = '  -2.41%  '
$ws.Range('E4').Value = '  -0.11%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '237.57'
$ws.Range('E5').Value = '  -2.05%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.607'
$ws.Range('E6').Value = '  -2.84%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '71.76'
$ws.Range('E7').Value = '  -2.95%  '
$ws.Range('E8').Value = '  -0.18%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.577'
$ws.Range('E9').Value = '  -4.04%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '39.66'
$ws.Range('E10').Value = '  -6.79%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0905'
$ws.Range('E11').Value = '  -5.27%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '54.43'
$ws.Range('E12').Value = '  -3.74%  '
$ws.Range('E13').Value = '  -3.02%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '6.69'
$ws.Range('E14').Value = '  -3.44%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '2.487.47'
$ws.Range('E15').Value = '  -2.48%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '14.27'
$ws.Range('E16').Value = '  +0.29%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '2.142.07'
$ws.Range('E17').Value = '  -3.42%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.776'
$ws.Range('E18').Value = '  -6.99%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '41.377.75'
$ws.Range('E19').Value = '  -1.22%  '
$ws.Range('E20').Value = '  -2.46%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '69.96'
$ws.Range('E21').Value = '  -3.87%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '5.76'
$ws.Range('E22').Value = '  -6.87%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '9.80'
$ws.Range('E23').Value = '  -11.40%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '226.10'
$ws.Range('E24').Value = '  -1.40%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.02'
$ws.Range('E25').Value = '  -2.92%  '
$ws.Range('E26').Value = '  +0.10%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '10.68'
$ws.Range('E27').Value = '  -5.95%  '
$ws.Range('E28').Value = '  -9.69%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.18'
$ws.Range('E29').Value = '  -3.72%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '2.15'
$ws.Range('E30').Value = '  -1.87%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '170.98'
$ws.Range('E31').Value = '  +2.47%  '
$ws.Range('B32').Value = 'EthereumClassic'
$ws.Range('C32').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '19.78'
$ws.Range('E32').Value = '  -3.52%  '
$ws.Range('B33').Value = 'InjectiveProtocol'
$ws.Range('C33').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '33.15'
$ws.Range('E33').Value = '  +11.01%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.0768'
$ws.Range('E34').Value = '  -3.67%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '5.27'
$ws.Range('E35').Value = '  -6.66%  '
$ws.Range('E36').Value = '  -3.34%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '4.28'
$ws.Range('E37').Value = '  -0.29%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.102'
$ws.Range('E38').Value = '  -6.23%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.0303'
$ws.Range('E39').Value = '  +0.49%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '12.02'
$ws.Range('E40').Value = '  -8.68%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '2.09'
$ws.Range('E41').Value = '  -1.46%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '5.35'
$ws.Range('E42').Value = '  -5.56%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '58.81'
$ws.Range('E43').Value = '  -8.95%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '8.46'
$ws.Range('E44').Value = '  -2.46%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.187'
$ws.Range('E45').Value = '  -5.24%  '
$ws.Range('E46').Value = '  -3.44%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '96.92'
$ws.Range('E47').Value = '  -6.56%  '
$ws.Range('E48').Value = '  -3.14%  '
$ws.Range('E49').Value = '  -4.79%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '2.17'
$ws.Range('E50').Value = '  -6.97%  '
$ws.Range('E51').Value = '  -2.33%  '
